$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 54378.105
$ws.Range("J70").Value = 2007.2307
$ws.Range("L70").Value = 6021.6921
$ws.Range("N70").Value = -6561.6921

$ws.Range("H73").Value = 54378.105
$ws.Range("J73").Value = 2007.2307
$ws.Range("L73").Value = 6021.6921
$ws.Range("N73").Value = -7893.6921

$ws.Range("H80").Value = 12613.223
$ws.Range("I80").Value = 400
$ws.Range("J80").Value = 14139.875
$ws.Range("K80").Value = 1200
$ws.Range("L80").Value = 42419.625
$ws.Range("M80").Value = -202
$ws.Range("N80").Value = -44415.625

$ws.Range("H83").Value = 12613.223
$ws.Range("I83").Value = 400
$ws.Range("J83").Value = 14139.875
$ws.Range("K83").Value = 3600
$ws.Range("L83").Value = 127258.875
$ws.Range("M83").Value = 1392
$ws.Range("N83").Value = -137242.875

$ws.Range("H112").Value = 4431.231
$ws.Range("J112").Value = 4431.231
$ws.Range("L112").Value = 13293.693
$ws.Range("N112").Value = -15509.693

$ws.Range("H129").Value = 1111.8243
$ws.Range("I129").Value = 310
$ws.Range("J129").Value = 1195.597
$ws.Range("K129").Value = 930
$ws.Range("L129").Value = 3586.791
$ws.Range("M129").Value = 4070
$ws.Range("N129").Value = -13586.791

$ws.Range("H137").Value = 2326.3333
$ws.Range("I137").Value = 1621.6923
$ws.Range("J137").Value = 2865.1765
$ws.Range("K137").Value = 4865.0769
$ws.Range("L137").Value = 8595.529500000001
$ws.Range("M137").Value = -2315.0769
$ws.Range("N137").Value = -13695.5295

$ws.Range("H138").Value = 3340.0874
$ws.Range("I138").Value = 1763.697
$ws.Range("J138").Value = 4446.915
$ws.Range("K138").Value = 5291.090999999999
$ws.Range("L138").Value = 13340.745
$ws.Range("M138").Value = -151.0909999999994
$ws.Range("N138").Value = -23620.745

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1000000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""

$ws.Range("H39").Value = 12071.429
$ws.Range("I39").Value = 12071.429
$ws.Range("K39").Value = 12071.429
$ws.Range("M39").Value = -11551.429

$ws.Range("H45").Value = 2072.8572
$ws.Range("I45").Value = 1978
$ws.Range("J45").Value = 2199.3333
$ws.Range("K45").Value = 1978
$ws.Range("L45").Value = 2199.3333
$ws.Range("M45").Value = -1601
$ws.Range("N45").Value = -2953.3333

$ws.Range("H88").Value = 2665.5557
$ws.Range("I88").Value = 2122.5
$ws.Range("J88").Value = 3100
$ws.Range("K88").Value = 2122.5
$ws.Range("L88").Value = 3100
$ws.Range("M88").Value = -1716.5
$ws.Range("N88").Value = -3912

$ws.Range("H91").Value = 2665.5557
$ws.Range("I91").Value = 2122.5
$ws.Range("J91").Value = 3100
$ws.Range("K91").Value = 2122.5
$ws.Range("L91").Value = 3100
$ws.Range("M91").Value = -718.5
$ws.Range("N91").Value = -5908

$ws.Range("H116").Value = 1000000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = ""

$ws.Range("H122").Value = 4958.9697
$ws.Range("I122").Value = 5598.923
$ws.Range("J122").Value = 2582
$ws.Range("K122").Value = 16796.769
$ws.Range("L122").Value = 7746
$ws.Range("M122").Value = -14346.769
$ws.Range("N122").Value = -12646

$ws.Range("H132").Value = 9190.471
$ws.Range("I132").Value = 10245.643
$ws.Range("J132").Value = 4266.3335
$ws.Range("K132").Value = 30736.929
$ws.Range("L132").Value = 12799.0005
$ws.Range("M132").Value = -28206.929
$ws.Range("N132").Value = -17859.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1000000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""

$ws.Range("H55").Value = 26140
$ws.Range("J55").Value = 26140
$ws.Range("L55").Value = 26140
$ws.Range("N55").Value = -26686

$ws.Range("H99").Value = 1717.1177
$ws.Range("I99").Value = 1448.3334
$ws.Range("J99").Value = 2362.2
$ws.Range("K99").Value = 1448.3334
$ws.Range("L99").Value = 2362.2
$ws.Range("M99").Value = 49.66660000000002
$ws.Range("N99").Value = -5358.2

$ws.Range("H107").Value = 3208.077
$ws.Range("I107").Value = 3373.182
$ws.Range("K107").Value = 3373.182
$ws.Range("M107").Value = -1453.182

$ws.Range("H134").Value = 2261.7368
$ws.Range("I134").Value = 2351.6428
$ws.Range("J134").Value = 2010
$ws.Range("K134").Value = 7054.928400000001
$ws.Range("L134").Value = 6030
$ws.Range("M134").Value = -4519.928400000001
$ws.Range("N134").Value = -11100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 747.7222
$ws.Range("I8").Value = 747.7222
$ws.Range("K8").Value = 2243.1666
$ws.Range("M8").Value = -2104.1666

$ws.Range("H103").Value = 2609.8572
$ws.Range("I103").Value = 2453.8
$ws.Range("K103").Value = 7361.400000000001
$ws.Range("M103").Value = -6482.400000000001

$ws.Range("H113").Value = 233250.1
$ws.Range("J113").Value = 812.65216
$ws.Range("L113").Value = 2437.95648
$ws.Range("N113").Value = -6777.95648

$ws.Range("H131").Value = 2663.2166
$ws.Range("I131").Value = 498.41177
$ws.Range("K131").Value = 1495.23531
$ws.Range("M131").Value = 3544.76469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3453.8948
$ws.Range("I132").Value = 2311.2
$ws.Range("J132").Value = 4723.5557
$ws.Range("K132").Value = 6933.599999999999
$ws.Range("L132").Value = 14170.6671
$ws.Range("M132").Value = -4403.599999999999
$ws.Range("N132").Value = -19230.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 28388.9
$ws.Range("J42").Value = 25861.125
$ws.Range("L42").Value = 25861.125
$ws.Range("N42").Value = -26987.125

$ws.Range("H49").Value = 28388.9
$ws.Range("J49").Value = 25861.125
$ws.Range("L49").Value = 25861.125
$ws.Range("N49").Value = -26155.125

$ws.Range("H69").Value = 131912.6
$ws.Range("J69").Value = 131912.6
$ws.Range("L69").Value = 131912.6
$ws.Range("N69").Value = -133534.6

$ws.Range("H72").Value = 131912.6
$ws.Range("J72").Value = 131912.6
$ws.Range("L72").Value = 395737.8
$ws.Range("N72").Value = -403849.8

$ws.Range("H132").Value = 3919.1904
$ws.Range("I132").Value = 3848.8647
$ws.Range("J132").Value = 4439.6
$ws.Range("K132").Value = 11546.5941
$ws.Range("L132").Value = 13318.8
$ws.Range("M132").Value = -9016.5941
$ws.Range("N132").Value = -18378.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 41682.332
$ws.Range("I58").Value = 50000
$ws.Range("K58").Value = 50000
$ws.Range("M58").Value = -49692

$ws.Range("H132").Value = 2996.6853
$ws.Range("I132").Value = 2663.3022
$ws.Range("J132").Value = 4299.909
$ws.Range("K132").Value = 7989.9066
$ws.Range("L132").Value = 12899.727
$ws.Range("M132").Value = -5459.9066
$ws.Range("N132").Value = -17959.727
